$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Name = "alpha2F"

$ws.Range("D13").Value = 0.9951749546634753
$ws.Range("J13").Value = 0.9951749546634753
$ws.Range("K13").Value = 0.995229341316833
